# Atualização da aula 5
# Move the "Imagem 13" picture on slide 4 upward (change its vertical
# position) - equivalent to changing <a:off y="2461450"/> to y="108545"
# (EMU) in the slide's XML, i.e. Top = 108545 / 12700 points.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(12)
$shape.Top = 8.546851
